# Project 1 code till Sep 5th
# Adds a new "TC006" (Select Hotel) test-data block (rows 10-11, columns A-I)
# to the "regression" sheet, mirroring the existing TC001/TC002/TC005 blocks.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("regression")

# --- Row 10: header row for TC006 -------------------------------------
# Column order chosen so new shared-strings are minted in the same order
# as the authored workbook (A, then E/E, F/F, G/H, G/H, I, I).
$ws.Range("A10").Value = "TC006"
$ws.Range("B10").Value = "username"
$ws.Range("C10").Value = "password"
$ws.Range("D10").Value = "expectedTitle"
$ws.Range("E10").Value = "Location"
$ws.Range("E11").Value = "London"
$ws.Range("F10").Value = "Number of Rooms"
$ws.Range("F11").Value = "2 - Two"
$ws.Range("G10").Value = "Check In Date"
$ws.Range("H10").Value = "Check Out Date"

# Dates entered as quote-prefixed text ('06/09/2024) so Excel stores them
# as literal strings (quotePrefix style) instead of date serials.
$ws.Range("G11").Value = "'06/09/2024"
$ws.Range("H11").Value = "'07/09/2024"

$ws.Range("I11").Value = "Adactin.com - Select Hotel"
$ws.Range("I10").Value = "expectedTitle2"

# --- Row 11: data row for TC006 ----------------------------------------
$ws.Range("A11").Value = "TC006"
$ws.Range("B11").Value = "reyaz0806"
$ws.Range("C11").Value = "reyaz123"
$ws.Range("D11").Value = "Adactin.com - Search Hotel"

# --- Column widths for the newly-used columns E:I -----------------------
$ws.Columns.Item(5).ColumnWidth = 16.053385416666668
$ws.Columns.Item(6).ColumnWidth = 34.276041666666664
$ws.Columns.Item(7).ColumnWidth = 25.830729166666668
$ws.Columns.Item(8).ColumnWidth = 29.166666666666668
$ws.Columns.Item(9).ColumnWidth = 49.053385416666664

# --- Selection moves to A10, matching the saved workbook view -----------
$ws.Range("A10").Select() | Out-Null

Write-Output "TC006 block added"
